$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.388.18"
$ws.Range("E2").Value = "  +4.29%  "
$ws.Range("D3").Value = "2.043.68"
$ws.Range("E3").Value = "  +2.85%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.650"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.91%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "65.76"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +9.85%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.401"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +10.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "59.37"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0812"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +9.59%  "
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.913"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.80"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +25.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.79"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.32%  "
$ws.Range("D16").Value = "2.344.59"
$ws.Range("E16").Value = "  +2.97%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.69"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.39%  "
$ws.Range("D18").Value = "2.046.66"
$ws.Range("E18").Value = "  +3.01%  "
$ws.Range("D19").Value = "37.333.65"
$ws.Range("E19").Value = "  +4.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "73.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.30%  "
$ws.Range("D21").Value = "0.0₃0901"
$ws.Range("E21").Value = "  +6.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.67%  "
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("E25").Value = "  +1.42%  "
$ws.Range("E26").Value = "  +4.65%  "
$ws.Range("E27").Value = "  +4.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.85%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.08"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.130"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +31.25%  "
$ws.Range("E31").Value = "  +3.00%  "
$ws.Range("E32").Value = "  +2.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0628"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.74%  "
$ws.Range("E35").Value = "  +5.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.40"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +12.57%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.36"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.40%  "
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("E39").Value = "  +3.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.04"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +31.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.30"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.12%  "
$ws.Range("E42").Value = "  +8.19%  "
$ws.Range("E43").Value = "  +5.36%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.16"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.06%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.44"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0219"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "95.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.82"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.82%  "
$ws.Range("D49").Value = "1.391.07"
$ws.Range("E49").Value = "  +2.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.94"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "47.09"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.87%  "
